# Update the "Prix Spot" worksheet with a new day column (AF) of EPEX spot prices,
# mirroring the style (header formatting) already used by the preceding date column (AE).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# New header label for column AF, row 1 - copy style from AE1 (bold header w/ border)
$ws.Range("AF1").Value = "15-jul"
$ws.Range("AE1").Copy()
$ws.Range("AF1").PasteSpecial(-4122) # xlPasteFormats

# New price values for column AF, rows 2-25
$values = @{
    2  = 76.77
    3  = 57.2
    4  = 63.01
    5  = 59.29
    6  = 54.1
    7  = 63.78
    8  = 69.52
    9  = 84.47
    10 = 94.16
    11 = 72.97
    12 = 30.5
    13 = 39.87
    14 = 46.2
    15 = 34.82
    16 = 25.8
    17 = 21.39
    18 = 20.82
    19 = 38.4
    20 = 44.77
    21 = 65.89
    22 = 74.45999999999999
    23 = 79.56999999999999
    24 = 100.79
    25 = 78.19
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 32).Value = $values[$row]
}
